# [Kadastro App] Yeni kayıt eklendi: 2882
# Adds a new record row (record No. 2882) to both the master "Kayitlar"
# sheet and the relevant "Erdemli" district sheet.

$wb = $excel.ActiveWorkbook

$recordNo   = "2882"
$recordDate = "2025-09-04"
$birim      = "Erdemli"
$parselSayi = "1"
$is         = "MAKS"
$personel   = "EMİNE ALANLI KIRCILI (K.Mühendisi), AYHAN KARADAYI (K.Teknisyeni)"

function Add-KadastroRow($SheetName) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Find the next empty row based on column A (skips the header row).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    # Every column in this table is stored as plain text (even the
    # numeric-looking / date-looking ones), matching the rest of the sheet.
    # Force a text number format first so Excel doesn't auto-coerce these
    # into numbers or dates.
    $rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 6))
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($newRow, 1).Value = $recordNo
    $ws.Cells.Item($newRow, 2).Value = $recordDate
    $ws.Cells.Item($newRow, 3).Value = $birim
    $ws.Cells.Item($newRow, 4).Value = $parselSayi
    $ws.Cells.Item($newRow, 5).Value = $is
    $ws.Cells.Item($newRow, 6).Value = $personel

    # Keep suppressing the "number stored as text" warning (xlNumberAsText)
    # over the whole used range, same as the rest of the table.
    $usedRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($newRow, 6))
    $usedRange.Errors.Item(3).Ignore = $true
}

Add-KadastroRow "Kayitlar"
Add-KadastroRow "Erdemli"
